# datasets with tweaked reward approach for boundary conditions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated notes text (shared string used by G4 & G5) ---
$ws.Range("G4").Value = 'Num2Move = 1; Discount of -1 for "far out" movements'
$ws.Range("G5").Value = 'Num2Move = 1; Discount of -1 for "far out" movements'

# --- Updated Iteration-3 / Iteration-4 metrics (row 4 & row 5) ---
$ws.Range("B4").Value = 219.84
$ws.Range("C4").Value = 19.940000000000001
$ws.Range("D4").Value = 6.74
$ws.Range("E4").Value = 8490.4
$ws.Range("F4").Value = 182.86

$ws.Range("B5").Value = 129.19
$ws.Range("C5").Value = 13.7
$ws.Range("D5").Value = 2.38
$ws.Range("E5").Value = 7740.5

# --- Column widths (approximate the re-sized / re-fitted columns) ---
$ws.Columns("A").ColumnWidth = 12.42578125
$ws.Columns("B").ColumnWidth = 12.08
$ws.Columns("C").ColumnWidth = 7.25
$ws.Columns("D").ColumnWidth = 12.58
$ws.Columns("E").ColumnWidth = 16.42
$ws.Columns("F").ColumnWidth = 13.92
$ws.Columns("G").ColumnWidth = 49.140625

# --- Header row: center everything horizontally, switch numeric columns to integer / one-decimal formats ---
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("B1").NumberFormat = "0"
$ws.Range("E1").NumberFormat = "0"
$ws.Range("F1").NumberFormat = "0"
$ws.Range("C1").NumberFormat = "0.0"
$ws.Range("D1").NumberFormat = "0.0"

# --- Data rows 2-5: center numeric columns & apply integer / one-decimal number formats ---
$ws.Range("B2:F5").HorizontalAlignment = -4108

foreach ($r in 2..5) {
    $ws.Cells.Item($r, 2).NumberFormat = "0"     # B - Num. Steps
    $ws.Cells.Item($r, 3).NumberFormat = "0.0"   # C - Error
    $ws.Cells.Item($r, 4).NumberFormat = "0.0"   # D - CEP Radius
    $ws.Cells.Item($r, 5).NumberFormat = "0"     # E - Distance Moved
    $ws.Cells.Item($r, 6).NumberFormat = "0"     # F - Total Reward
}

# --- Row 4 lost its italic "needs review" styling, now the normal font ---
$ws.Range("A4").Font.Italic = $false

# --- Row 5 is now fully bold+italic (finalized / highlighted metrics row) ---
$ws.Range("A5:G5").Font.Bold = $true
$ws.Range("A5:G5").Font.Italic = $true

# --- Conditional highlight fills: green = good boundary, amber = needs attention ---
$green = 5296274   # RGB(146, 208, 80) -> FF92D050
$amber = 49407      # RGB(255, 192, 0)  -> FFFFC000

$ws.Range("B2:C2").Interior.Color = $green
$ws.Range("E3:F3").Interior.Color = $amber
$ws.Range("B4:D4").Interior.Color = $amber
$ws.Range("F4").Interior.Color = $green
$ws.Range("D5:E5").Interior.Color = $green

# --- Selection cursor moved ---
$ws.Range("J3").Select()

Write-Host "edit applied"
